$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rename the header row: "<Name>_old" -> "<Name>_FV2410" (cols A-J) and
#    "<Name>_new" -> "<Name>_FV2504" (cols L-U). Column K ("diff") unchanged.
# ---------------------------------------------------------------------------
$baseNames = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $oldCol = $i + 1          # columns 1..10  -> A..J
    $newCol = $i + 12         # columns 12..21 -> L..U
    $ws.Cells.Item(1, $oldCol).Value = ($baseNames[$i] + "_FV2410")
    $ws.Cells.Item(1, $newCol).Value = ($baseNames[$i] + "_FV2504")
}

# ---------------------------------------------------------------------------
# 2. Freeze the header row (split at row 1, top-left cell A2).
# ---------------------------------------------------------------------------
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# ---------------------------------------------------------------------------
# 3. Turn the used range A1:U52 into an Excel Table ("Table1") so that a
#    xl/tables/table1.xml part + tableParts reference is emitted.
# ---------------------------------------------------------------------------
$rng = $ws.Range("A1:U52")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $rng, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

$ws.Range("A1").Select() | Out-Null
